$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Fill in new "RF" (faltas recorrentes?) column J values for each student row
$ws.Range("J2").Value = 0
$ws.Range("J3").Value = 1
$ws.Range("J4").Value = 1
$ws.Range("J5").Value = 0
$ws.Range("J6").Value = 1

# Mark row 5 (Jordan Santos Hohenfeld) with "RF" flag in column N, like row 2
$ws.Range("N5").Value = "RF"

# Update the active selection to J3 as in the saved file
$ws.Range("J3").Select()
